# Weekly update: insert two new daily price rows for "Poroto verde" at the
# top of the data block (rows 635:636), shifting the existing rows down by
# two. Excel's Insert() shifts the rows below and copies formatting (e.g.
# the date-number-format style) from the row above the insertion point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 635:636 - everything currently at 635:689 moves
# down to 637:691.
$ws.Rows("635:636").Insert()

# ---- New row 635 ----
$ws.Cells.Item(635, 1).Value = 6
$ws.Cells.Item(635, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(635, 3).Value = "Metropolitana"
$ws.Cells.Item(635, 4).Value = 44578
$ws.Cells.Item(635, 5).Value = 13
$ws.Cells.Item(635, 6).Value = 100112031
$ws.Cells.Item(635, 7).Value = "Poroto verde"
$ws.Cells.Item(635, 8).Value = "Magnum"
$ws.Cells.Item(635, 9).Value = "Primera"
$ws.Cells.Item(635, 10).Value = 270
$ws.Cells.Item(635, 11).Value = 25000
$ws.Cells.Item(635, 12).Value = 27000
$ws.Cells.Item(635, 13).Value = 26111
$ws.Cells.Item(635, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(635, 15).Value = "Región Metropolitana"
$ws.Cells.Item(635, 16).Value = 1044
$ws.Cells.Item(635, 17).Value = 25
$ws.Cells.Item(635, 18).Value = "Hortaliza"

# ---- New row 636 ----
$ws.Cells.Item(636, 1).Value = 6
$ws.Cells.Item(636, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(636, 3).Value = "Metropolitana"
$ws.Cells.Item(636, 4).Value = 44578
$ws.Cells.Item(636, 5).Value = 13
$ws.Cells.Item(636, 6).Value = 100112031
$ws.Cells.Item(636, 7).Value = "Poroto verde"
$ws.Cells.Item(636, 8).Value = "Sin especificar"
$ws.Cells.Item(636, 9).Value = "Primera"
$ws.Cells.Item(636, 10).Value = 770
$ws.Cells.Item(636, 11).Value = 35000
$ws.Cells.Item(636, 12).Value = 37000
$ws.Cells.Item(636, 13).Value = 36013
$ws.Cells.Item(636, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(636, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(636, 16).Value = 1441
$ws.Cells.Item(636, 17).Value = 25
$ws.Cells.Item(636, 18).Value = "Hortaliza"
